$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{
        B = -20.60755986554534
        C = 1.089226581658092
        D = 75.00302793421463
        E = -0.406779944259445
        F = -0.352054911594391
        G = 0.7815408847023604
        H = -0.4500861918092851
        I = 0.8988809271980212
        J = 4.3552525155391
        K = 27
        L = -13.08276063154555
        N = 4.35525251498856
        O = 5.021919181655227
    }
    3 = @{
        B = 0.000006489415692342842
        C = 0.7035650081437954
        D = 1.195271404766907
        E = 62.69770586818203
        F = 1.653471968311442
        G = -1.727242461161577
        H = 0.5345901828952315
        I = -0.6299647063015685
        J = 4.355350724842692
        K = 98
        L = -10.91721351340661
        N = 4.355350721035633
        O = 5.0220173877023
    }
    4 = @{
        B = -1.836450284879223
        C = 62.8767700094226
        D = 0.1799469150252427
        E = -0.0001415292567473507
        F = -0.6433970046019519
        G = -0.7616181374878477
        H = 0.8604019339048006
        I = 1.848922505396178
        J = 4.355362752407631
        K = 17
        L = -3.972176334657421
        N = 4.35536275239051
        O = 5.022029419057177
    }
    5 = @{
        B = 0.3725194491593553
        C = 8.517556335889296
        D = 58.72967226374491
        E = -0.0290272403757472
        F = 0.7927274845000563
        G = -1.324838884042288
        H = -0.7077263909663822
        I = 1.073934080018848
        J = 4.355367582331654
        K = 77
        L = -5.976831165322857
        N = 4.355367582335208
        O = 5.022034249001875
    }
    6 = @{
        B = -0.3630679222453754
        C = 1.444344719566855
        D = 47.84348870804013
        E = 20.08156574257011
        F = 0.8703964381228366
        G = 0.7040910487935337
        H = -0.3947333488908169
        I = -0.8932278747492548
        J = 4.354054096912488
        K = 40
        L = -18.49099562563747
        N = 4.355476733141792
        O = 5.022143399808459
    }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
